$d = $word.ActiveDocument

# 1. Highlight "Perfis dos observados" in green
$range = $d.Content
$found = $range.Find.Execute("Perfis dos observados", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $range.Font.HighlightColorIndex = 4
}

# 2. Highlight "Crie uma person" (part of "Crie uma persona") in green.
#    This splits the original run into "Crie uma person" (highlighted) and
#    "a para cada perfil a ser observado" (not highlighted).
$range2 = $d.Content
$found2 = $range2.Find.Execute("Crie uma person", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $range2.Font.HighlightColorIndex = 4
}

# 3. Remove the existing _GoBack bookmark (previously located after
#    "...a maioria do conteúdo que recebe.")
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 4. Re-create the _GoBack bookmark between "Crie uma persona" and
#    " para cada perfil a ser observado" (i.e. right after "persona").
$range3 = $d.Content
$found3 = $range3.Find.Execute("Crie uma persona", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found3) {
    $insertPoint = $range3.End
    $bmRange = $d.Range($insertPoint, $insertPoint)
    $d.Bookmarks.Add("_GoBack", $bmRange)
}
